$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 0.3299849534320388
$ws.Range("B2").Value = 0.4184978746830859
$ws.Range("C2").Value = 0.2417807408181229
$ws.Range("D2").Value = 0.3737467453250149
$ws.Range("E2").Value = 0.2839621223878922
$ws.Range("K2").Value = 2.309894674024271
$ws.Range("L2").Value = 2.929485122781602
$ws.Range("M2").Value = 1.69246518572686
$ws.Range("N2").Value = 2.616227217275104
$ws.Range("O2").Value = 1.987734856715246
$ws.Range("P2").Value = 66.14035
$ws.Range("Q2").Value = 108.8290731272669
$ws.Range("R2").Value = 43.95103301855618
$ws.Range("S2").Value = 79.48393468757878
$ws.Range("T2").Value = 51.63600791664881
$ws.Range("U2").Value = 0.2019516717667799
$ws.Range("V2").Value = 0.2922138194345176
$ws.Range("W2").Value = 0.09985910553342195
$ws.Range("X2").Value = 0.2508071157700391
$ws.Range("Y2").Value = 0.1478193981814509
$ws.Range("Z2").Value = 0.8501516043801202
$ws.Range("AA2").Value = 0.9446876151351921
$ws.Range("AB2").Value = 0.6855584441312089
$ws.Range("AC2").Value = 0.9011788995946686
$ws.Range("AD2").Value = 0.7794325914978461
$ws.Range("F3").Value = 6.994794719503335
$ws.Range("G3").Value = 8.854542141129887
$ws.Range("H3").Value = 5.132826858615592
$ws.Range("I3").Value = 7.923351523289907
$ws.Range("J3").Value = 6.028623807981445
$ws.Range("K3").Value = 2.3082822574361
$ws.Range("L3").Value = 2.921998906572862
$ws.Range("M3").Value = 1.693832863343145
$ws.Range("N3").Value = 2.614706002685669
$ws.Range("O3").Value = 1.989445856633877
$ws.Range("P3").Value = 64.17094
$ws.Range("Q3").Value = 80.85807097720003
$ws.Range("R3").Value = 56.18609427057385
$ws.Range("S3").Value = 68.87462406001431
$ws.Range("T3").Value = 58.84640390942783
$ws.Range("U3").Value = 0.2017544353779299
$ws.Range("V3").Value = 0.2917376639182558
$ws.Range("W3").Value = 0.1000069737949224
$ws.Range("X3").Value = 0.2503411223966134
$ws.Range("Y3").Value = 0.147910713417528
$ws.Range("Z3").Value = 0.8504914315656196
$ws.Range("AA3").Value = 0.943211181054956
$ws.Range("AB3").Value = 0.6920824332762201
$ws.Range("AC3").Value = 0.9015301455825644
$ws.Range("AD3").Value = 0.7809126469286877
$ws.Range("A4").Value = 0.3301105416098413
$ws.Range("B4").Value = 0.4179068582842576
$ws.Range("C4").Value = 0.2416849141139535
$ws.Range("D4").Value = 0.3740157998440143
$ws.Range("E4").Value = 0.2843093394167646
$ws.Range("F4").Value = 7.000603011245695
$ws.Range("G4").Value = 8.858955327440841
$ws.Range("H4").Value = 5.125928625277665
$ws.Range("I4").Value = 7.923112155530301
$ws.Range("J4").Value = 6.034413403095476
$ws.Range("K4").Value = 2.31095479881165
$ws.Range("L4").Value = 3.241848093887591
$ws.Range("M4").Value = 1.495681560499509
$ws.Range("N4").Value = 2.746066621702917
$ws.Range("O4").Value = 1.858766723696207
$ws.Range("P4").Value = 67.36566000000001
$ws.Range("Q4").Value = 116.774986235109
$ws.Range("R4").Value = 42.87138240793082
$ws.Range("S4").Value = 82.03059320246304
$ws.Range("T4").Value = 51.32892910210004
$ws.Range("U4").Value = 0.1997598320828628
$ws.Range("V4").Value = 0.3285547090205936
$ws.Range("W4").Value = 0.06362921624793981
$ws.Range("X4").Value = 0.2693643256374609
$ws.Range("Y4").Value = 0.1249662934311058
$ws.Range("Z4").Value = 0.8357140874367625
$ws.Range("AA4").Value = 0.9723777907646807
$ws.Range("AB4").Value = 0.5721374368607285
$ws.Range("AC4").Value = 0.9087312320516334
$ws.Range("AD4").Value = 0.7279159123833915
